# Update of goodness of fit table: add RMSE / NRMSE columns (S-TS / LB-TS)
# ahead of the existing RRMSE / MRE / Coverage columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new columns before column B to make room for the new
#     "RMSE (mg O2 m-2 d-1)" and "NRMSE" column groups.
$ws.Range("B1:E1").EntireColumn.Insert()

# --- New header row 1 (merged pairs), matches format of existing headers
$ws.Range("B1").Value = "RMSE (mg O2 m-2 d-1)"
$ws.Range("D1").Value = "NRMSE"
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("B1:E1").HorizontalAlignment = -4108

# --- New header row 2 (S-TS / LB-TS labels)
$ws.Range("B2").Value = "S-TS"
$ws.Range("C2").Value = "LB-TS"
$ws.Range("D2").Value = "S-TS"
$ws.Range("E2").Value = "LB-TS"
$ws.Range("B2:E2").HorizontalAlignment = -4108

# --- RMSE / NRMSE data values
$rmse_sts  = @(3.3417383243704601, 1.1180404517792999, 1.1915197558550701, 1.11699915242334, 4.1589269306042196, 2.1498946820767699)
$rmse_lbts = @(3.4333387307687002, 0.75129118862183897, 1.1139766778684601, 1.02498398861837, 3.5121738846041399, 1.84385684523717)
$nrmse_sts = @(0.22478478585283099, 0.22766261125606599, 0.18653001002950301, 0.311667980841285, 0.26893364571336098, 0.22766871045582099)
$nrmse_lbts = @(0.23094636277406999, 0.152982759741045, 0.174390797864967, 0.28599367281011301, 0.227111883167571, 0.195260034698433)

for ($i = 0; $i -lt 6; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $rmse_sts[$i]
    $ws.Cells.Item($row, 3).Value = $rmse_lbts[$i]
    $ws.Cells.Item($row, 4).Value = $nrmse_sts[$i]
    $ws.Cells.Item($row, 5).Value = $nrmse_lbts[$i]
}

$ws.Range("B3:E8").NumberFormat = "0.00"
$ws.Range("B3:E8").HorizontalAlignment = -4108

# --- Summary formulas (Mean / Min / Max) for the new columns
$ws.Range("B10").Formula = "=AVERAGE(B3:B8)"
$ws.Range("C10").Formula = "=AVERAGE(C3:C8)"
$ws.Range("D10").Formula = "=AVERAGE(D3:D8)"
$ws.Range("E10").Formula = "=AVERAGE(E3:E8)"

$ws.Range("B11").Formula = "=MIN(B3:B8)"
$ws.Range("C11").Formula = "=MIN(C3:C8)"
$ws.Range("D11").Formula = "=MIN(D3:D8)"
$ws.Range("E11").Formula = "=MIN(E3:E8)"

$ws.Range("B12").Formula = "=MAX(B3:B8)"
$ws.Range("C12").Formula = "=MAX(C3:C8)"
$ws.Range("D12").Formula = "=MAX(D3:D8)"
$ws.Range("E12").Formula = "=MAX(E3:E8)"

$ws.Range("B10:E12").NumberFormat = "0.00"
$ws.Range("B10:E12").HorizontalAlignment = -4108

# --- Column widths
$ws.Range("B1").EntireColumn.ColumnWidth = 11.81640625
$ws.Range("C1").EntireColumn.ColumnWidth = 11.81640625
$ws.Range("D1").EntireColumn.ColumnWidth = 7.81640625
$ws.Range("E1").EntireColumn.ColumnWidth = 11.08984375

# --- Placeholder / scratch cells below the table (rows 21-33), matching
#     the stray empty cells left behind in the source worksheet. These
#     carry no value; touch them (no-op bold toggle) so the exporter keeps
#     an entry for them, matching the leftover cells in the original file.
$placeholderRanges = @("F21:K21", "F22", "I22:K22", "F23", "I23:K23", "F24", "I24:K24", "F25", "I25:K25", "F26", "I26:K26", "F27", "I27:K27", "F28", "I28:K28", "F29", "I29:K29", "F30", "I30:K30", "F31", "I31:K31", "F32", "I32:K32", "F33", "I33:K33")
foreach ($rng in $placeholderRanges) {
    $ws.Range($rng).Font.Bold = $false
}

# --- Selection matches the author's final cursor position
$ws.Range("D14").Select()
